$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D.  This shifts the existing
# NUM_OF_PAGES / MESSAGE_TO_CONNECTION columns one place to the right
# (D->E, E->F) and makes room for the new "STARTING_PAGE" column.
$ws.Columns.Item(4).Insert()

# The inserted column inherits the formatting (wrapped header style) of
# the old column D - drop it back to the plain/default style used by the
# rest of the data cells.
$ws.Columns.Item(4).ClearFormats()

# New header + value for scenario 2's starting page.
$ws.Range("D1").Value = "STARTING_PAGE"
$ws.Range("D2").Value = 2

# New header for scenario 3's message body column (shifted into F).
$ws.Range("F1").Value = "MESSAGE_TO_CONNECTION"

# Widen the new last column so its header text fits the row without
# wrapping (keeps row 1 at the default height).
$ws.Columns.Item(6).ColumnWidth = 28.5

# The header row no longer needs the extra wrapped-text height now that
# the wide text fits on one line - shrink it back down to the default.
$ws.Rows.Item(1).AutoFit()

$ws.Range("D2").Select()

$wb.Save()
